$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Info")

# Insert a new row 1 on the Info sheet (pushes Type/Number/DocDate/StartDate/DeliveryDate down by one row)
$ws.Rows("1:1").Insert()

# Give the new row 1 the same look (bold label / highlighted value) as the row below it
$ws.Range("A2:B2").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)

# Populate the new "Title" row
$ws.Range("A1").Value = "Title"
$ws.Range("B1").Value = "Document Title"

# Widen column B slightly to fit the new, longer value
$ws.Columns("B").ColumnWidth = 14.5703125

# Keep the workbook-level defined names pointing at the right cells after the insert
$wb.Names.Item("Info!DeliveryDate").RefersTo = "=Info!`$B`$6"
$wb.Names.Item("Info!DocDate").RefersTo = "=Info!`$B`$4"
$wb.Names.Item("Info!Number").RefersTo = "=Info!`$B`$3"
$wb.Names.Item("Info!StartDate").RefersTo = "=Info!`$B`$5"
$wb.Names.Item("Info!Type").RefersTo = "=Info!`$B`$2"
$ws.Names.Add("Title", "=Info!`$B`$1")

# Make the Info sheet the active tab, with B1 selected
$ws.Activate()
$ws.Range("B1").Select()
